# plot_curve_with_activation.xlsx - "Add files via upload" edit
#
# The author re-uploaded the workbook with:
#  - a second "lotes_7" / new "lotes_8" pair of header labels added to the
#    right of the existing table header row (columns AC:AH), plus some
#    blank " " spacer header cells inserted in row 1 (F1, J1, S1, AD1, AH1)
#  - a batch of "second factor" input cells (the *_2 column of each
#    lotes_n pair) zeroed out for several rows, which ripples through the
#    adjacent SUM() formulas
#  - the view/selection state left pointing at a different cell
#  - column X given an explicit width

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new header cells in row 1 -------------------------------------------
# (space spacer columns reuse the existing shared string " ")
$ws.Range("F1").Value = " "
$ws.Range("J1").Value = " "
$ws.Range("S1").Value = " "
$ws.Range("AD1").Value = " "
$ws.Range("AH1").Value = " "
# new label columns - AC1 repeats the existing "lotes_7" label,
# AG1 introduces the brand new "lotes_8" label/shared string
$ws.Range("AC1").Value = "lotes_7"
$ws.Range("AG1").Value = "lotes_8"

# --- zero out the "second" input cell of several lotes_n pairs -----------
# Column D (lotes_1 pair), rows 4,6,7,8,9,10
$ws.Range("D4").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("D10").Value = 0

# Column G (lotes_2 pair), rows 9,10
$ws.Range("G9").Value = 0
$ws.Range("G10").Value = 0

# Column M (lotes_3 pair), row 10
$ws.Range("M10").Value = 0

# Column P (lotes_4 pair), rows 9,10
$ws.Range("P9").Value = 0
$ws.Range("P10").Value = 0

# Column S (lotes_5 pair), rows 9,10
$ws.Range("S9").Value = 0
$ws.Range("S10").Value = 0

# Column V (lotes_6 pair), row 10
$ws.Range("V10").Value = 0

# --- column width for the newly-relevant column X -------------------------
$ws.Columns("X").ColumnWidth = 10.65

# --- selection / view state -----------------------------------------------
$ws.Range("S13").Select() | Out-Null
